$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (Q) that mirrors the formatting of the existing
# "2019" column (P), then fill in the new data values.
$ws.Range("P3:P8").Copy()
$ws.Range("Q3:Q8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("Q3").Value = $null
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 25.6
$ws.Range("Q6").Value = 13.073527219449954
$ws.Range("Q7").Value = 21.941290626870046
$ws.Range("Q8").Value = 196.6

# Row 1 gets a taller, explicit row height in the updated sheet.
$ws.Rows.Item(1).RowHeight = 19.5
